# Update "data pobrania" (download timestamp) column V for every data row
# on the sheet to reflect the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-01-06 10:47:35"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 22).Value = $newTimestamp
}
